# Update Name of Algo
# Applies the numeric value changes described by the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.564499999999995
$ws.Range("A10").Value = -22.0771
$ws.Range("A12").Value = -22.05759999999999
$ws.Range("B12").Value = 5.985100000000001
$ws.Range("C12").Value = -12.50480000000001
$ws.Range("C13").Value = -12.3871
$ws.Range("B17").Value = 5.165699999999996
$ws.Range("A18").Value = -22.1469
$ws.Range("C21").Value = -13.7677
$ws.Range("B26").Value = 4.438400000000002
$ws.Range("B27").Value = 6.475100000000003
$ws.Range("B28").Value = 6.252199999999997
$ws.Range("C36").Value = -11.65260000000001
$ws.Range("A37").Value = -21.92189999999999
$ws.Range("B37").Value = 6.5291
$ws.Range("C38").Value = -12.2422
$ws.Range("C41").Value = -12.96070000000002
$ws.Range("C52").Value = -10.9566
$ws.Range("A55").Value = -22.0362
$ws.Range("B65").Value = 5.7761
$ws.Range("C67").Value = -11.3507
$ws.Range("A68").Value = -21.4663
$ws.Range("B73").Value = 9.397599999999997
$ws.Range("A77").Value = -20.7496
$ws.Range("A78").Value = -19.97699999999999
$ws.Range("B84").Value = 5.526200000000001
$ws.Range("B85").Value = 5.325
$ws.Range("C89").Value = -14.33079999999999
$ws.Range("B93").Value = 5.605500000000003
$ws.Range("B95").Value = 6.1454
$ws.Range("C95").Value = -12.6208
$ws.Range("B98").Value = 4.943400000000009
$ws.Range("B99").Value = 5.513199999999999
$ws.Range("B101").Value = 5.845999999999998
$ws.Range("C105").Value = -12.56880000000001
